# Update cryptos list with latest scraped price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "29.049.69"
$ws.Range("E2").Value2 = "  -0.41%  "

$ws.Range("D3").Value2 = "1.825.21"
$ws.Range("E3").Value2 = "  -0.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value2 = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.09"
$ws.Range("E5").Value2 = "  -1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6148"
$ws.Range("E6").Value2 = "  -2.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value2 = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07343"
$ws.Range("E8").Value2 = "  -2.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2883"
$ws.Range("E9").Value2 = "  -1.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.78"
$ws.Range("E10").Value2 = "  -1.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07703"
$ws.Range("E11").Value2 = "  -0.62%  "

$ws.Range("D12").Value2 = "1.806.00"
$ws.Range("E12").Value2 = "  -1.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.924"
$ws.Range("E13").Value2 = "  -1.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6601"
$ws.Range("E14").Value2 = "  -1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "81.71"
$ws.Range("E15").Value2 = "  -1.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008958"
$ws.Range("E16").Value2 = "  -4.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.839"
$ws.Range("E17").Value2 = "  -2.93%  "

$ws.Range("D18").Value2 = "29.045.70"
$ws.Range("E18").Value2 = "  -0.39%  "

$ws.Range("D19").Value2 = "2.066.10"
$ws.Range("E19").Value2 = "  -0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.85"
$ws.Range("E20").Value2 = "  +5.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("E21").Value2 = "  -1.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.007"
$ws.Range("E22").Value2 = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.083"
$ws.Range("E23").Value2 = "  -0.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.008"
$ws.Range("E24").Value2 = "  +0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.81"
$ws.Range("E25").Value2 = "  -1.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1401"
$ws.Range("E26").Value2 = "  +0.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.418"
$ws.Range("E27").Value2 = "  -1.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.61"
$ws.Range("E28").Value2 = "  -2.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.489"
$ws.Range("E29").Value2 = "  -1.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05563"
$ws.Range("E30").Value2 = "  -5.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.077"
$ws.Range("E31").Value2 = "  +0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.078"
$ws.Range("E32").Value2 = "  -1.95%  "

$ws.Range("E33").Value2 = "  +0.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.829"
$ws.Range("E34").Value2 = "  -1.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7319"
$ws.Range("E35").Value2 = "  -2.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.132"
$ws.Range("E36").Value2 = "  -0.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.640"
$ws.Range("E37").Value2 = "  -1.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.817"
$ws.Range("E38").Value2 = "  +1.64%  "

$ws.Range("B39").Value2 = "Maker"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value2 = "1.198.65"
$ws.Range("E39").Value2 = "  -2.44%  "

$ws.Range("B40").Value2 = "VeChain"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01755"
$ws.Range("E40").Value2 = "  -2.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.379"
$ws.Range("E41").Value2 = "  -2.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8825"
$ws.Range("E42").Value2 = "  -1.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.006"
$ws.Range("E43").Value2 = "  -0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.00"
$ws.Range("E44").Value2 = "  -1.32%  "

$ws.Range("D45").Value2 = "1.957.99"
$ws.Range("E45").Value2 = "  -1.01%  "

$ws.Range("B46").Value2 = "Aave"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.15"
$ws.Range("E46").Value2 = "  -2.50%  "

$ws.Range("B47").Value2 = "Mantle"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5112"
$ws.Range("E47").Value2 = "  -0.03%  "

$ws.Range("E48").Value2 = "  -3.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3984"
$ws.Range("E49").Value2 = "  -2.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.028"
$ws.Range("E50").Value2 = "  +0.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05797"
$ws.Range("E51").Value2 = "  -0.86%  "
